$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.824.90'
$ws.Range("E2").Value = '  +1.40%  '

# Row 3
$ws.Range("D3").Value = '2.048.46'
$ws.Range("E3").Value = '  +0.98%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '''229.58'
$ws.Range("E5").Value = '  +0.99%  '

# Row 6
$ws.Range("D6").Value = '''0.614'
$ws.Range("E6").Value = '  +0.67%  '

# Row 7
$ws.Range("D7").Value = '''58.33'
$ws.Range("E7").Value = '  +5.42%  '

# Row 8
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("E9").Value = '  +1.51%  '

# Row 10
$ws.Range("E10").Value = '  +2.04%  '

# Row 11
$ws.Range("E11").Value = '  +0.57%  '

# Row 12
$ws.Range("D12").Value = '2.350.94'
$ws.Range("E12").Value = '  +1.02%  '

# Row 13
$ws.Range("D13").Value = '''14.55'
$ws.Range("E13").Value = '  +1.96%  '

# Row 14
$ws.Range("D14").Value = '''20.66'
$ws.Range("E14").Value = '  +1.29%  '

# Row 15
$ws.Range("E15").Value = '  +1.75%  '

# Row 16
$ws.Range("E16").Value = '  +0.38%  '

# Row 17
$ws.Range("D17").Value = '2.048.43'
$ws.Range("E17").Value = '  +1.21%  '

# Row 18
$ws.Range("D18").Value = '37.770.33'
$ws.Range("E18").Value = '  +1.45%  '

# Row 19
$ws.Range("E19").Value = '  -2.32%  '

# Row 20
$ws.Range("D20").Value = '''69.67'

# Row 21
$ws.Range("E21").Value = '  +1.38%  '

# Row 22
$ws.Range("D22").Value = '''224.13'
$ws.Range("E22").Value = '  -0.18%  '

# Row 23
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
$ws.Range("E24").Value = '  +0.26%  '

# Row 25
$ws.Range("E25").Value = '  +2.26%  '

# Row 26
$ws.Range("D26").Value = '''166.68'
$ws.Range("E26").Value = '  +0.87%  '

# Row 27
$ws.Range("E27").Value = '  -0.94%  '

# Row 28
$ws.Range("E28").Value = '  +2.88%  '

# Row 29
$ws.Range("D29").Value = '''18.97'
$ws.Range("E29").Value = '  +0.76%  '

# Row 30
$ws.Range("E30").Value = '  -0.92%  '

# Row 31
$ws.Range("E31").Value = '  +1.23%  '

# Row 32
$ws.Range("E32").Value = '  -0.20%  '

# Row 33
$ws.Range("D33").Value = '''2.08'
$ws.Range("E33").Value = '  +13.10%  '

# Row 34
$ws.Range("E34").Value = '  +2.43%  '

# Row 35
$ws.Range("D35").Value = '''0.0611'
$ws.Range("E35").Value = '  -1.21%  '

# Row 36
$ws.Range("E36").Value = '  -1.21%  '

# Row 37
$ws.Range("E37").Value = '  +9.64%  '

# Row 38
$ws.Range("D38").Value = '''3.30'
$ws.Range("E38").Value = '  +4.10%  '

# Row 39
$ws.Range("D39").Value = '''0.999'
$ws.Range("E39").Value = '  -0.08%  '

# Row 40
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.484.89'
$ws.Range("E40").Value = '  +0.47%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '''0.0217'
$ws.Range("E41").Value = '  -0.04%  '

# Row 42
$ws.Range("D42").Value = '''97.20'
$ws.Range("E42").Value = '  +0.90%  '

# Row 43
$ws.Range("E43").Value = '  +2.64%  '

# Row 44
$ws.Range("D44").Value = '''0.0931'
$ws.Range("E44").Value = '  +0.86%  '

# Row 45
$ws.Range("E45").Value = '  -0.32%  '

# Row 46
$ws.Range("D46").Value = '''4.18'
$ws.Range("E46").Value = '  +15.92%  '

# Row 48
$ws.Range("E48").Value = '  -0.46%  '

# Row 49
$ws.Range("D49").Value = '''2.96'
$ws.Range("E49").Value = '  +0.97%  '

# Row 50
$ws.Range("E50").Value = '  -3.62%  '

# Row 51
$ws.Range("D51").Value = '2.239.01'
$ws.Range("E51").Value = '  +1.07%  '
